$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 存款 (deposits) sheet

# ---------------------------------------------------------------------------
# 1. Apply the existing cell styles to the new columns (G:M) before filling
#    in values, so formatting matches the rest of the sheet (header style on
#    row 1, data-row style on rows 2-8).
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("G2:M8").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2. Row 1 becomes a genuine header row (field names), matching the
#    convention used on the other sheets, instead of duplicating row 2's
#    data values.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# ---------------------------------------------------------------------------
# 3. Column I (date) must stay a literal text value "2012-04-30" rather than
#    be auto-converted into a date serial number, so force text format on
#    that range before typing the values in.
# ---------------------------------------------------------------------------
$ws.Range("I2:I8").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 4. Fill in the new metadata columns for every data row (2-8): the
#    property_category/category/date/legislator_name/legislator_id/
#    source_file are constant for this whole import, while index mirrors
#    column A of the same row.
# ---------------------------------------------------------------------------
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 7).Value  = "deposit"                     # G: property_category
    $ws.Cells.Item($r, 8).Value  = "normal"                      # H: category
    $ws.Cells.Item($r, 9).Value  = "2012-04-30"                  # I: date
    $ws.Cells.Item($r, 10).Value = "蘇震清"                       # J: legislator_name
    $ws.Cells.Item($r, 11).Value = 1718                          # K: legislator_id
    $ws.Cells.Item($r, 12).Value = "tmp16a71"                    # L: source_file
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($r, 1).Value2  # M: index (= column A)
}

# Re-apply the data-row style over column I so the temporary text number
# format doesn't linger as a one-off style distinct from the rest of the row.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("I2:I8").PasteSpecial(-4122) | Out-Null
